# VS26 Isidore Quantum Quad Chart Text - refine to concise/scannable format
# matching the quad-chart template (per commit message).
#
# The whole slide-text block lives as a single run (uniform 20pt formatting)
# inside the document's 3rd paragraph, with `w:br` line breaks (character
# 0x0B / backtick-v in PowerShell) separating each line. We replace that
# paragraph's full text in one shot so Word re-splits it into the correct
# <w:t>/<w:br/> run sequence automatically, preserving the run formatting.

$d = $word.ActiveDocument
$p3 = $d.Paragraphs(3).Range

$newText = @"
VALIANT SHIELD 2026 EXPERIMENTATION PACKET`vIsidore Quantum® Maritime Cybersecurity Solution`v`v===== SLIDE 1: EXPERIMENT OVERVIEW =====`v`vEXPERIMENT TITLE:`v(U) Quantum-Safe Maritime Cybersecurity - Isidore Quantum®`v`vDESCRIPTION:`vValidate NSA-engineered post-quantum encryption platform for maritime operations. Addresses dual threats: current cyber attacks (spoofing, ransomware, malware) and future quantum-enabled adversaries. Drop-in integration with zero infrastructure retrofits.`v`vOBJECTIVE:`vDemonstrate TRL 8-9 advancement in operational naval environment with autonomous threat detection and CNSA 2.0 compliance.`v`vAVAILABILITY DATE:`vQ2-Q3 2026 (Exercise execution 28 JUN 2026)`v`vPROGRAM RISKS / IMPACT:`vMinimal - drop-in design, autonomous post-deployment. Early validation supports Navy modernization priorities.`v`vGEOGRAPHICAL DOMAIN:`vMaritime (vessels, ports, satellite links)`v`vRECOMMENDED HOST UNIT:`vPACFLT (Third Fleet or Seventh Fleet operational assets)`v`vEXPERIMENTATION SPONSOR:`vUnited States Pacific Fleet (PACFLT). Supporting: USINDOPACOM J81 JEESC, Naval Supply Systems Command.`v`vFUNDING STATUS:`vForward Edge-AI funded. NSA partnership supports R&D costs.`v`vPOC INFORMATION:`vBrandon@forwardedge.ai | forwardedge.ai`v`v`v===== SLIDE 2: DEPLOYMENT & REQUIREMENTS =====`v`vDEPLOYMENT:`v2+ operational naval vessels (bridge/AIS systems), 1 major WESTPAC port facility, 1 satellite gateway. Forward Edge-AI technical team for integration support. Duration: 4-6 weeks during VS26 exercise window.`v`vHOST COMPONENT / UNIT SUPPORT REQUESTED:`vVessel bridge access and systems integration. Port logistics platform access. Satellite ground station coordination. Observer seats for 6-8 personnel (cybersecurity SMEs, command staff, logistics personnel). Power and network integration support.`v`vEXERCISE CONSIDERATIONS:`vEncrypted threat data handling (classified). Autonomous operation minimizes host unit burden post-integration. No infrastructure retrofits required. Real-time threat logs require secure reporting infrastructure.`v`vREQUIRED EXERCISE EVENTS:`vOperational vessel underway time; port cargo operations; satellite communication exercises.`v`vPERSONNEL / EQUIPMENT LIST:`vForward Edge-AI integration team (2-3)`vNaval IT personnel (4-5)`vPort operations staff (2-3)`vObservers (6-8)`vIsidore Quantum units (variants: Enterprise, Standard, IoT)`vPower: <12W per unit`vFootprint: 140x89x39mm`vLatency: <90µs
"@

$p3.Text = $newText
Write-Output "Paragraph 3 replaced. New length: $($p3.Text.Length)"
